# "Scroll added in log off"
# The ManageProducts seed sheet has its randomly generated SKU placeholders
# (column B, rows 2/5/8) refreshed with new values, same as the on-disk
# "thin top+bottom border / white fill" look the cells already had.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-ProductCell($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)

    # Re-apply the same visual style (thin border top & bottom, white fill)
    # that these generated SKU cells already carry.
    $rng.Borders.Item(8).LineStyle = 1
    $rng.Borders.Item(8).Weight = 2
    $rng.Borders.Item(9).LineStyle = 1
    $rng.Borders.Item(9).Weight = 2
    $rng.Interior.ColorIndex = 2

    $rng.Value = $newValue
}

Set-ProductCell "B2" "prodisjY"
Set-ProductCell "B5" "prodZedD"
Set-ProductCell "B8" "prodtJAD"
